$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 6)
$ws.Range("D2").Value = 44223
$ws.Range("K2").Value = 'Ruby Diamond'
$ws.Range("R2").Value = 'Región Metropolitana'

# Row 3 (was row 7)
$ws.Range("D3").Value = 44223
$ws.Range("K3").Value = 'Super Queen'
$ws.Range("L3").Value = 'Segunda'

# Row 4 (was row 3)
$ws.Range("D4").Value = 44243
$ws.Range("K4").Value = 'Venus'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 250
$ws.Range("R4").Value = 'Región de O''Higgins'

# Row 5 (was row 26)
$ws.Range("D5").Value = 44215
$ws.Range("K5").Value = 'Venus'
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'

# Row 6 (was row 33)
$ws.Range("D6").Value = 44273
$ws.Range("K6").Value = 'Artic Snow'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 23000
$ws.Range("P6").Value = 22500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 1250

# Row 7 (was row 34)
$ws.Range("D7").Value = 44273
$ws.Range("K7").Value = 'August Red'
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("S7").Value = 1250

# Row 8 (was row 27)
$ws.Range("D8").Value = 44168
$ws.Range("K8").Value = 'Artic Star'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 270
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("S8").Value = 1306

# Row 9 (was row 28)
$ws.Range("D9").Value = 44168
$ws.Range("K9").Value = 'Early Glo'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 23000
$ws.Range("O9").Value = 24000
$ws.Range("P9").Value = 23500
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1306

# Row 10 (was row 29)
$ws.Range("D10").Value = 44217
$ws.Range("K10").Value = 'Ruby Diamond'
$ws.Range("Q10").Value = '$/caja 18 kilos empedrada'

# Row 11 (was row 30)
$ws.Range("D11").Value = 44217
$ws.Range("K11").Value = 'Venus'
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18500
$ws.Range("S11").Value = 1028

# Row 12 (was row 31)
$ws.Range("D12").Value = 44244
$ws.Range("K12").Value = 'Nectar Crest'
$ws.Range("Q12").Value = '$/caja 18 kilos granel'

# Row 13 (was row 32)
$ws.Range("D13").Value = 44244
$ws.Range("K13").Value = 'Venus'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("S13").Value = 1083

# Row 14 (was row 17)
$ws.Range("D14").Value = 44202
$ws.Range("K14").Value = 'Super Queen'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21000
$ws.Range("S14").Value = 1167

# Row 15 (was row 9)
$ws.Range("D15").Value = 44167
$ws.Range("K15").Value = 'Early John'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 25000
$ws.Range("O15").Value = 26000
$ws.Range("P15").Value = 25500
$ws.Range("R15").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S15").Value = 1417

# Row 16 (was row 4)
$ws.Range("D16").Value = 44174
$ws.Range("K16").Value = 'Early John'
$ws.Range("M16").Value = 200
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Provincia de San Felipe de Aconcagua'

# Row 17 (was row 22)
$ws.Range("D17").Value = 44278
$ws.Range("K17").Value = 'August Red'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 23500
$ws.Range("Q17").Value = '$/bandeja 18 kilos granel'
$ws.Range("S17").Value = 1306

# Row 18 (was row 23)
$ws.Range("D18").Value = 44278
$ws.Range("K18").Value = 'June Pearl'
$ws.Range("L18").Value = 'Primera'
$ws.Range("N18").Value = 23000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 23500
$ws.Range("S18").Value = 1306

# Row 19 (was row 8)
$ws.Range("D19").Value = 44209
$ws.Range("K19").Value = 'Super Queen'
$ws.Range("L19").Value = 'Tercera'
$ws.Range("M19").Value = 320
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 17500
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("S19").Value = 972

# Row 20 (was row 24)
$ws.Range("D20").Value = 44229
$ws.Range("K20").Value = 'Artic Sprite'
$ws.Range("M20").Value = 300

# Row 21 (was row 20)
$ws.Range("D21").Value = 44216
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("S21").Value = 1083

# Row 22 (was row 5)
$ws.Range("D22").Value = 44257
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 19000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19500
$ws.Range("Q22").Value = '$/caja 18 kilos granel'
$ws.Range("S22").Value = 1083

# Row 23 (was row 18)
$ws.Range("D23").Value = 44169
$ws.Range("K23").Value = 'Artic Sprite'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 25000
$ws.Range("P23").Value = 24500
$ws.Range("S23").Value = 1361

# Row 24 (was row 19)
$ws.Range("D24").Value = 44169
$ws.Range("K24").Value = 'Early John'
$ws.Range("M24").Value = 270
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 24500
$ws.Range("S24").Value = 1361

# Row 25 (was row 21)
$ws.Range("D25").Value = 44222
$ws.Range("K25").Value = 'Nectar Crest'
$ws.Range("M25").Value = 270
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20500
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("S25").Value = 1139

# Row 26 (was row 10)
$ws.Range("D26").Value = 44210
$ws.Range("K26").Value = 'Early John'

# Row 27 (was row 11)
$ws.Range("D27").Value = 44210
$ws.Range("K27").Value = 'Nectar Crest'
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 19000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 19500
$ws.Range("Q27").Value = '$/bandeja 18 kilos granel'
$ws.Range("S27").Value = 1083

# Row 28 (was row 12)
$ws.Range("D28").Value = 44210
$ws.Range("K28").Value = 'Red Diamond'
$ws.Range("N28").Value = 19000
$ws.Range("O28").Value = 20000
$ws.Range("P28").Value = 19500
$ws.Range("Q28").Value = '$/bandeja 18 kilos granel'
$ws.Range("S28").Value = 1083

# Row 29 (was row 13)
$ws.Range("D29").Value = 44161
$ws.Range("K29").Value = 'Artic Glo'
$ws.Range("M29").Value = 280
$ws.Range("N29").Value = 25000
$ws.Range("O29").Value = 26000
$ws.Range("P29").Value = 25500
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("S29").Value = 1417

# Row 30 (was row 14)
$ws.Range("D30").Value = 44161
$ws.Range("K30").Value = 'Early John'
$ws.Range("N30").Value = 26000
$ws.Range("O30").Value = 27000
$ws.Range("P30").Value = 26500
$ws.Range("Q30").Value = '$/caja 18 kilos granel'
$ws.Range("S30").Value = 1472

# Row 31 (was row 15)
$ws.Range("D31").Value = 44238
$ws.Range("K31").Value = 'August Red'
$ws.Range("M31").Value = 320
$ws.Range("N31").Value = 20000
$ws.Range("O31").Value = 21000
$ws.Range("P31").Value = 20500
$ws.Range("S31").Value = 1139

# Row 32 (was row 16)
$ws.Range("D32").Value = 44238
$ws.Range("M32").Value = 320
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("Q32").Value = '$/bandeja 18 kilos granel'
$ws.Range("S32").Value = 1139

# Row 33 (was row 2)
$ws.Range("D33").Value = 44236
$ws.Range("K33").Value = 'June Pearl'
$ws.Range("M33").Value = 270
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 21000
$ws.Range("P33").Value = 20500
$ws.Range("Q33").Value = '$/caja 18 kilos granel'
$ws.Range("S33").Value = 1139

# Row 34 (was row 25)
$ws.Range("D34").Value = 44201
$ws.Range("K34").Value = 'Super Queen'
$ws.Range("Q34").Value = '$/caja 18 kilos granel'
